{"js": "// The worksheet of division equations lives in the single table in the\n// document body. Only the five \"content\" rows (1, 5, 9, 13, 17 - 1-based)\n// hold text; the others are blank spacer rows. Every populated cell's\n// equation is replaced with a new one, in reading order (row major,\n// left-to-right), while leaving the run/paragraph formatting untouched.\nconst newValues = [\n  [\"50\u00f76=8, 2\", \"54\u00f74=13, 2\", \"92\u00f76=15, 2\", \"32\u00f78=4, 0\", \"50\u00f72=25, 0\"],\n  [\"76\u00f72=38, 0\", \"11\u00f76=1, 5\", \"75\u00f74=18, 3\", \"28\u00f79=3, 1\", \"15\u00f74=3, 3\"],\n  [\"41\u00f75=8, 1\", \"74\u00f73=24, 2\", \"92\u00f78=11, 4\", \"56\u00f72=28, 0\", \"79\u00f72=39, 1\"],\n  [\"19\u00f79=2, 1\", \"62\u00f78=7, 6\", \"82\u00f79=9, 1\", \"33\u00f76=5, 3\", \"84\u00f74=21, 0\"],\n  [\"98\u00f72=49, 0\", \"27\u00f74=6, 3\", \"51\u00f79=5, 6\", \"50\u00f74=12, 2\", \"71\u00f73=23, 2\"],\n];\n\nconst rowIndexes = [0, 4, 8, 12, 16]; // 0-based table-row indexes that hold data\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (let r = 0; r < rowIndexes.length; r++) {\n  const tableRow = rowIndexes[r];\n  for (let c = 0; c < newValues[r].length; c++) {\n    const cell = table.getCell(tableRow, c);\n    const range = cell.body.getRange();\n    range.insertText(newValues[r][c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The worksheet of division equations lives in the single table in the\n# document. Only five rows (1, 5, 9, 13, 17) hold text; the rest are blank\n# spacer rows. Every populated cell's equation is replaced with a new one,\n# in reading order (row major, left-to-right), while run/paragraph\n# formatting is left untouched (Range.Text only replaces the text run).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rows = @(1, 5, 9, 13, 17)\n$newValues = @(\n    @(\"50\u00f76=8, 2\", \"54\u00f74=13, 2\", \"92\u00f76=15, 2\", \"32\u00f78=4, 0\", \"50\u00f72=25, 0\"),\n    @(\"76\u00f72=38, 0\", \"11\u00f76=1, 5\", \"75\u00f74=18, 3\", \"28\u00f79=3, 1\", \"15\u00f74=3, 3\"),\n    @(\"41\u00f75=8, 1\", \"74\u00f73=24, 2\", \"92\u00f78=11, 4\", \"56\u00f72=28, 0\", \"79\u00f72=39, 1\"),\n    @(\"19\u00f79=2, 1\", \"62\u00f78=7, 6\", \"82\u00f79=9, 1\", \"33\u00f76=5, 3\", \"84\u00f74=21, 0\"),\n    @(\"98\u00f72=49, 0\", \"27\u00f74=6, 3\", \"51\u00f79=5, 6\", \"50\u00f74=12, 2\", \"71\u00f73=23, 2\")\n)\n\nfor ($r = 0; $r -lt $rows.Length; $r++) {\n    $row = $rows[$r]\n    $values = $newValues[$r]\n    for ($c = 1; $c -le $values.Length; $c++) {\n        $cell = $t.Cell($row, $c)\n        $cell.Range.Text = $values[$c - 1]\n    }\n}\n"}
